$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.332.74"
$ws.Range("E2").Value = "  -3.30%  "

$ws.Range("D3").Value = "3.514.85"
$ws.Range("E3").Value = "  -4.83%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.98"
$ws.Range("E5").Value = "  -1.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.96"
$ws.Range("E6").Value = "  -3.80%  "

$ws.Range("E7").Value = "  +0.72%  "

$ws.Range("D8").Value = "3.506.48"
$ws.Range("E8").Value = "  -4.73%  "

$ws.Range("E9").Value = "  +0.04%  "

$ws.Range("E10").Value = "  -6.24%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.77"
$ws.Range("E11").Value = "  +5.05%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.596"
$ws.Range("E12").Value = "  -2.98%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.01"
$ws.Range("E13").Value = "  -6.02%  "

$ws.Range("E14").Value = "  -3.86%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "678.02"
$ws.Range("E15").Value = "  -0.94%  "

$ws.Range("D16").Value = "4.075.87"
$ws.Range("E16").Value = "  -4.99%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.75"
$ws.Range("E17").Value = "  -3.28%  "

$ws.Range("D18").Value = "69.281.37"
$ws.Range("E18").Value = "  -3.53%  "

$ws.Range("D19").Value = "3.513.47"
$ws.Range("E19").Value = "  -4.87%  "

$ws.Range("E20").Value = "  -1.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.48"
$ws.Range("E21").Value = "  -3.85%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.21"
$ws.Range("E22").Value = "  -4.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.905"
$ws.Range("E23").Value = "  -4.37%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.15"
$ws.Range("E24").Value = "  -9.52%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.09"
$ws.Range("E25").Value = "  -5.62%  "

$ws.Range("E26").Value = "  -4.44%  "

$ws.Range("E27").Value = "  -0.67%  "

$ws.Range("E28").Value = "  +0.12%  "

$ws.Range("E29").Value = "  -6.53%  "

$ws.Range("E30").Value = "  -7.91%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.05"
$ws.Range("E31").Value = "  -6.62%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.74"
$ws.Range("E32").Value = "  -5.94%  "

$ws.Range("E33").Value = "  -7.97%  "

$ws.Range("E34").Value = "  -1.10%  "

$ws.Range("E35").Value = "  -6.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "595.53"
$ws.Range("E36").Value = "  +5.31%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.61"
$ws.Range("E37").Value = "  -15.57%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.92"
$ws.Range("E38").Value = "  -3.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.105"
$ws.Range("E39").Value = "  -4.66%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "57.26"
$ws.Range("E40").Value = "  -3.85%  "

$ws.Range("E41").Value = "  +0.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0440"
$ws.Range("E42").Value = "  -5.91%  "

$ws.Range("E43").Value = "  -4.69%  "

$ws.Range("E44").Value = "  -6.32%  "

$ws.Range("D45").Value = "3.415.84"
$ws.Range("E45").Value = "  -9.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "33.38"
$ws.Range("E46").Value = "  -6.55%  "

$ws.Range("D47").Value = "0.0₃0709"
$ws.Range("E47").Value = "  -8.78%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.92"
$ws.Range("E48").Value = "  +0.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.61"
$ws.Range("E49").Value = "  -7.23%  "

$ws.Range("E50").Value = "  -0.70%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.80"
$ws.Range("E51").Value = "  +18.32%  "
